$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 ("비추천수" = dislike count): copy E1's header style (bold, centered,
# bordered), then set the text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "비추천수"

# Data cells F2:F91 hold the dislike counts. They mirror the existing columns
# (e.g. E = 추천수/like count): numeric-looking values stored as text, so force
# the Text number format before writing them in, one row per comment.
$ws.Range("F2:F91").NumberFormat = "@"

$dislikes = @("4","1","1","3","0","2","1","1","1","1","1","0","1","0","1","1","1","1","0","1","0","0","0","0","0","0","1","1","1","1","0","0","0","0","1","1","1","1","1","1","0","0","0","0","0","0","0","0","0","0","0","0","0","0","1","1","1","1","1","1","1","1","1","0","0","0","0","0","0","0","0","0","0","0","0","0","2","1","1","1","1","1","1","1","1","1","1","1","1","2")

for ($i = 0; $i -lt $dislikes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $dislikes[$i]
}
